$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the daily date blocks forward by one week:
# 2023-09-18..23 -> 2023-09-25..30 (each date spans 10 rows of tasks)
$ws.Range("A2:A11").Value = "2023-09-25"
$ws.Range("A12:A21").Value = "2023-09-26"
$ws.Range("A22:A31").Value = "2023-09-27"
$ws.Range("A32:A41").Value = "2023-09-28"
$ws.Range("A42:A51").Value = "2023-09-29"
$ws.Range("A52:A61").Value = "2023-09-30"

# Update the saved view/selection to match the new scroll position
$win = $excel.ActiveWindow
$win.ScrollRow = 43
$win.ScrollColumn = 1
$ws.Range("B52").Select()
